$d = $word.ActiveDocument

# This document prints one merged record (MERGEFIELD results already
# resolved to literal text) laid out twice - once per table cell on the
# page. The edit moves the printed record forward (NO 5 -> 2) and updates
# the corresponding name / id / size values, leaving every other field
# (GENDER, KELAS, UB_2, UB_4, UB_5, ...) untouched.
#
# Each field's Result is updated in-place via a Find scoped strictly to
# that field's own Result range (not the whole document) so a short
# numeric match (e.g. "5") can never bleed into a neighbouring field
# whose value happens to contain that digit (e.g. "25"). Using
# wdReplaceOne (1) instead of wdReplaceAll (2) keeps the single
# in-range match from cascading into a second, unintended match.

$replacements = @{
    "MERGEFIELD NO"   = @{ old = "5";              new = "2" }
    "MERGEFIELD NAMA" = @{ old = "FATHUL BARI";     new = "SANDY SATRIA WIDJAYA" }
    "MERGEFIELD TOPI" = @{ old = "2020.01.2.0024";  new = "2020.01.2.0008" }
    "MERGEFIELD UB_1" = @{ old = "44";              new = "43" }
    "MERGEFIELD UB_3" = @{ old = "16";              new = "15" }
    "MERGEFIELD UB_6" = @{ old = "26";              new = "25" }
    "MERGEFIELD UB_7" = @{ old = "71";              new = "70" }
    "MERGEFIELD UB_8" = @{ old = "38";              new = "36" }
}

foreach ($f in $d.Fields) {
    $code = $f.Code.Text.Trim()
    if ($replacements.ContainsKey($code)) {
        $entry = $replacements[$code]
        $resultRange = $f.Result
        $resultRange.Find.Execute($entry.old, $true, $false, $false, $false, $false, $true, 1, $false, $entry.new, 1)
    }
}
